$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'The game is an intergalactic murder mystery, which emulates a board game scenario in which players spin a wheel to determine a randomly selected number of moves by which they travel through space to various planets.'
$ws.Range("B2").Value = 'Each board is made up of cells'
$ws.Range("C2").Value = 'No, the two statements are not semantically equal. Statement 1 talks about a board being made up of cells, which seems to refer to a simple structural detail. Statement 2, however, describes a specific type of game and its gameplay, which is more complex and detailed. The two statements discuss different aspects of possibly different board games and are therefore not semantically equal.'
$ws.Range("D2").Value = 'No, these two statements do not convey the exact same information. The first statement speaks generally about a board being made up of cells. The second statement describes a specific type of game scenario, with much more detail provided about the game mechanics and theme.'
$ws.Range("E2").Value = 'No, these two statements are not conveying the same meaning. The first statement is describing the composition of a board, suggesting that it is divided into different sections or pieces, referred to as cells. The second statement is describing a specific type of game, with its various rules and characteristics. It does include the mention of a ''board game scenario'', but does not specify how the ''board'' in that game is composed or whether it includes ''cells''.'
$ws.Range("F2").Value = 'No, these two statements are not synonymous. The first statement is a general statement about boards being made up of cells. The second statement specifically describes a game and has no mention of the board being made up of cells.'
$ws.Range("G2").Value = 'No, statement 1 and statement 2 do not have identical implications. Statement 1 suggests that each board (possibly referring to a board game) has cells, but does not specify anything about the nature of the game or its mechanics. On the other hand, Statement 2 describes a specific type of game with a particular playing method, involving intergalactic travel and murder mystery, but it doesn''t indicate anything about the board being made up of cells.'

$ws.Range("A3").Value = 'These entry points are graphically depicted on the game as a diagonal line off a square on the board leading into a planet''s atmosphere.'
$ws.Range("B3").Value = 'Each board is made up of cells'
$ws.Range("C3").Value = 'No, the two statements are not semantically equal. The first statement describes the composition of a board, while the second statement describes how entry points in a game are depicted.'
$ws.Range("D3").Value = 'No, these two statements do not convey the exact same information. The first statement speaks about the composition of a board, while the second statement describes how certain graphical elements are represented in a game.'
$ws.Range("E3").Value = 'No, these two statements are not conveying the same meaning. The first statement is describing the composition of a board (being made up of cells), while the second statement is describing how entry points are represented in a game.'
$ws.Range("F3").Value = 'No, these statements are not synonymous. They are conveying two different pieces of information. Statement 1 is explaining the composition of the board whereas Statement 2 is describing how something is visually represented in a game.'
$ws.Range("G3").Value = 'No, statement 1 and statement 2 do not have identical implications. Statement 1 discusses the composition of a board, presumably in a game context, being made up of cells. Statement 2, however, is specific to a graphical element of a game relating to entry points depicted on a game board. They discuss different aspects in a gaming context and, therefore, don''t have identical implications.'

$ws.Range("A4").Value = 'A player visits a planet by making legal moves on the board until a cell adjacent to an atmosphere entry point is reached.'
$ws.Range("B4").Value = 'Each board is made up of cells'
$ws.Range("C4").Value = 'No, the two statements are not semantically equal. The first statement is talking about the composition of a board, while the second statement is describing a gameplay action within the context of a game.'
$ws.Range("D4").Value = 'No, the two statements do not convey the exact same information. The first statement describes what a board is made of, while the second statement describes a way a player interacts with the board in a game.'
$ws.Range("E4").Value = 'No, these two statements are not conveying the same meaning. The first statement is describing the composition of the board, while the second statement is providing instructions or details about a gameplay.'
$ws.Range("F4").Value = 'No, these statements are not synonymous. The first statement describes the components of a board, while the second one explains the process of a game or a game strategy.'
$ws.Range("G4").Value = 'No, statement 1 and statement 2 do not have identical implications. Statement 1 discusses the composition of the board, while statement 2 provides information about a game action.'

$ws.Range("A5").Value = 'If an announced hypothesis is incorrect, the player loses the game and cannot pose hypotheses any longer or make moves on the board but must continue refuting the hypotheses of other players.'
$ws.Range("B5").Value = 'Each board is made up of cells'
$ws.Range("C5").Value = 'No, the two statements are not semantically equal. The first statement describes the composition of a board, while the second statement describes the consequences of a particular action in a game. There is no semantic overlap or equality between these two statements.'
$ws.Range("D5").Value = 'No, these two statements do not convey the same information. The first statement talks about the composition of a board, while the second statement discusses what happens if a player in a game makes an incorrect hypothesis.'
$ws.Range("E5").Value = 'No, these two statements are not conveying the same meaning. The first statement is talking about the composition of a board, implying it''s formed from multiple cells. The second statement is outlining the consequences for a player in a game if their announced hypothesis is incorrect. They don''t convey the same idea or information.'
$ws.Range("F5").Value = 'No, these statements are not synonymous. They describe completely different scenarios. The first talks about the composition of a board, and the second outlines the consequences of a wrong hypothesis in a game.'
$ws.Range("G5").Value = 'No, statement 1 and statement 2 do not have identical implications. Statement 1 is about the composition of a board, while statement 2 is about the consequence of an incorrect hypothesis for a player in a game. They refer to completely different things and hence, their implications are different.'

$ws.Range("A6").Value = 'Each planet has between one and three predefined atmospheric entry points.'
$ws.Range("B6").Value = 'Each planet is made up of one to three planets'
$ws.Range("C6").Value = 'No, these two statements are not semantically equal. The first statement talks about the composition of a planet, mentioning a questionable fact that a planet is made of other planets. The second statement, on the other hand, discusses the number of predefined entry points into a planet''s atmosphere. The subjects are similar (both discussing planets), but the predicates - what''s being said about the subjects - are entirely different.'
$ws.Range("D6").Value = 'No, the two statements do not convey the same information. The first statement is about the composition of planets, while the second statement is about the number of atmospheric entry points a planet might have.'
$ws.Range("E6").Value = 'No, these two statements are not conveying the same meaning. The first statement is about composition of the planets, whereas the second statement is regarding the entry points in the atmosphere of the planets.'
$ws.Range("F6").Value = 'No, these statements are not synonymous. The first statement is about the composition of a planet implying it can somehow be made up of other planets, which scientifically incorrect. The second statement is about the number of predefined atmospheric entry points a planet has, which refers to space travel and exploration.'
$ws.Range("G6").Value = 'No, statement 1 and statement 2 do not have identical implications. Statement 1 suggests that each planet consists of one to three other planets, which is logically impossible. On the other hand, statement 2 implies that each planet has predefined points of atmospheric entry, which is a plausible claim related to space exploration. The two statements discuss completely different aspects of planets.'
